# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook's detail table (rows 16-34) listed debt periods for four
# workers (Noel Enrique Casseres Cassiani, Yulisa Castilla Lopez, Richard
# Rafael Hernandez Otaiza and Jorge Felix Hidalgo Marrugo). The new data
# drops Noel and Yulisa, keeps Jorge Felix Hidalgo Marrugo's single row,
# and replaces Richard's 16 months of arrears (previously listed newest
# period first, with two different "Valor Mora" amounts and a smaller
# "Salario Basico") with the same 16 periods listed oldest period first
# and refreshed amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two rows that belonged to workers who are no longer part of
# this "estado de cuenta" (Noel Enrique Casseres Cassiani and Yulisa
# Castilla Lopez). Everything below shifts up by two rows automatically,
# which is also what moves the closing signature block from rows 39/40
# to rows 37/38.
$ws.Rows("16:17").Delete()

# Row 16: Jorge Felix Hidalgo Marrugo keeps his original single period.
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "73572504"
$ws.Cells.Item(16, 4).Value = "JORGE FELIX HIDALGO MARRUGO"
$ws.Cells.Item(16, 5).Value = "2010"
$ws.Cells.Item(16, 6).Value = 1170
$ws.Cells.Item(16, 7).Value = 1500000

# Rows 17-32: Richard Rafael Hernandez Otaiza's 16 periods, now listed
# oldest-to-newest (2206 .. 2309) with the refreshed mora/salario values.
$richard = @(
  @("2206", 54845, 1510000),
  @("2207", 54845, 1510000),
  @("2208", 54845, 1510000),
  @("2209", 54845, 1510000),
  @("2210", 54845, 1510000),
  @("2211", 54845, 1510000),
  @("2212", 54845, 1510000),
  @("2301", 54845, 1510000),
  @("2302", 54845, 1510000),
  @("2303", 54845, 1510000),
  @("2304", 54845, 1510000),
  @("2305", 54845, 1510000),
  @("2306", 54845, 1510000),
  @("2307", 54845, 1510000),
  @("2308", 53976, 1510000),
  @("2309", 53976, 1510000)
)

$r = 17
foreach ($period in $richard) {
  $ws.Cells.Item($r, 2).Value = "PPT"
  $ws.Cells.Item($r, 3).Value = "5069162"
  $ws.Cells.Item($r, 4).Value = "RICHARD RAFAEL HERNANDEZ OTAIZA"
  $ws.Cells.Item($r, 5).Value = $period[0]
  $ws.Cells.Item($r, 6).Value = $period[1]
  $ws.Cells.Item($r, 7).Value = $period[2]
  $r = $r + 1
}

# Summary header: total mora, worker count and period count now reflect
# the trimmed-down table (2 workers / 17 rows / 16 periods for Richard).
$ws.Range("E11").Value = 876952
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 17
